# Auto-update predictions and index for 2025-10-24
#
# Two fixtures finished without a recorded prediction outcome and are
# dropped from the sheet (Shakhtar Donetsk - Legia Warszawa, AS Roma - FC
# Viktoria Plzen). The remaining fixtures are completed, so their rows are
# refreshed with final scores, updated AI confidence and a checkmark in
# the Result column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Shakhtar Donetsk - Legia Warszawa" row (row 2).
$ws.Rows.Item(2).Delete()

# After the shift-up, "AS Roma - FC Viktoria Plzen" is now row 3.
$ws.Rows.Item(3).Delete()

# Remaining rows (now 2-4): Maccabi Tel Aviv, Feyenoord Rotterdam, AEK Athens.
# Row 2: Maccabi Tel Aviv - FC Midtjylland finished 0:3.
$ws.Range("A2").Value = "Maccabi Tel Aviv - FC Midtjylland ✓: 0:3"
$ws.Range("C2").Value = 64
$ws.Range("G2").Value = "✓"

# Row 3: Feyenoord Rotterdam - Panathinaikos FC finished 3:1.
$ws.Range("A3").Value = "Feyenoord Rotterdam ✓ - Panathinaikos FC: 3:1"
$ws.Range("C3").Value = 63
$ws.Range("G3").Value = "✓"

# Row 4: AEK Athens - Aberdeen FC finished 6:0.
$ws.Range("A4").Value = "AEK Athens ✓ - Aberdeen FC: 6:0"
$ws.Range("C4").Value = 56
$ws.Range("G4").Value = "✓"
